# Correct the Excel import-test stub:
#  - rename header "ID_CATEGORIA" -> "CATEGORIA"
#  - E2 sample value becomes text "Categoria test" (was numeric 1)
#  - H2 sample value becomes text "11/01/2020" (was a date serial 44146),
#    while keeping its original mm/dd/yy-style numeric format applied
#  - selection moves to H3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename (stays a string, same cell / position).
$ws.Range("E1").Value = "CATEGORIA"

# Sample row: category becomes descriptive text instead of a raw id.
$ws.Range("E2").Value = "Categoria test"

# Sample row: publication date stored as literal text "11/01/2020"
# instead of a real date serial. Force text entry (so Excel doesn't
# reinterpret the string as a date), then restore the original
# mm/dd/yy number format on the cell.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "11/01/2020"
$ws.Range("H2").NumberFormat = "mm/dd/yy"

# Active cell moves to H3.
$ws.Range("H3").Select()
